$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aguilar Prototype")

# Helper: write a date-like string into a cell while preventing Excel's
# automatic "looks like a date" conversion, and without disturbing the
# destination cell's existing style. We stage the text in a scratch cell
# far outside the used range, format that scratch cell as Text so the
# string isn't reinterpreted, then copy/paste-special *values only* into
# the real destination (which keeps the destination's original style).
function Set-TextDate {
    param([string]$CellAddr, [string]$DateText)

    $scratch = $ws.Range("AA1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $DateText
    $scratch.Copy()
    $ws.Range($CellAddr).PasteSpecial(-4163) # xlPasteValues
    $excel.CutCopyMode = $false
    $scratch.Clear()
}

# Row 29: date update and Present value change only
Set-TextDate "N29" "2025-10-20"
$ws.Range("Q29").Value = 2.24

# Row 30: date update and full shift of Present/Lag1-4
Set-TextDate "N30" "2025-10-20"
$ws.Range("Q30").Value = 2.26
$ws.Range("R30").Value = 2.27
$ws.Range("S30").Value = 2.28
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.3

# Row 39: date update and full shift, U39 becomes empty
Set-TextDate "N39" "2025-10-17"
$ws.Range("Q39").Value = 121.1218
$ws.Range("R39").Value = 121.0834
$ws.Range("S39").Value = 121.2669
$ws.Range("T39").Value = 121.5815
$ws.Range("U39").ClearContents()

# Row 40: date update and Present value change only
Set-TextDate "N40" "2025-10-17"
$ws.Range("Q40").Value = -3.121703151457274

# Row 48: date update and full shift, T48 becomes populated, U48 becomes empty
Set-TextDate "N48" "2025-10-17"
$ws.Range("Q48").Value = 3.46
$ws.Range("R48").Value = 3.41
$ws.Range("S48").Value = 3.5
$ws.Range("T48").Value = 3.48
$ws.Range("U48").ClearContents()

# Row 49
Set-TextDate "N49" "2025-10-17"
$ws.Range("Q49").Value = 3.59
$ws.Range("R49").Value = 3.55
$ws.Range("S49").Value = 3.63
$ws.Range("T49").Value = 3.6
$ws.Range("U49").ClearContents()

# Row 50
Set-TextDate "N50" "2025-10-17"
$ws.Range("Q50").Value = 4.02
$ws.Range("R50").Value = 3.99
$ws.Range("S50").Value = 4.05
$ws.Range("T50").Value = 4.03
$ws.Range("U50").ClearContents()

# Row 52
Set-TextDate "N52" "2025-10-17"
$ws.Range("Q52").Value = 5.72
$ws.Range("R52").Value = 5.7
$ws.Range("S52").Value = 5.73
$ws.Range("T52").Value = 5.74
$ws.Range("U52").ClearContents()

$wb.Save()
